$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Markov transition-matrix probabilities after simulating more games.
$updates = @{
    "B2" = 0.1720116618075802
    "C2" = 0.5772594752186589
    "J2" = 0.02915451895043732
    "P2" = 0.1166180758017493
    "S2" = 0.1049562682215743
    "B3" = 0.004901960784313725
    "C3" = 0.009803921568627451
    "J3" = 0.06862745098039216
    "P3" = 0.75
    "S3" = 0.1666666666666667
    "J4" = 0.02222222222222222
    "P4" = 0.6666666666666666
    "S4" = 0.3111111111111111
    "B6" = 0.06403940886699508
    "D6" = 0.01477832512315271
    "F6" = 0.0541871921182266
    "J6" = 0.1970443349753695
    "O6" = 0.02463054187192118
    "Q6" = 0.2216748768472906
    "R6" = 0.06896551724137931
    "S6" = 0.354679802955665
    "B7" = 0.09195402298850575
    "D7" = 0.03448275862068965
    "E7" = 0.005747126436781609
    "F7" = 0.06896551724137931
    "J7" = 0.1551724137931035
    "O7" = 0.02298850574712644
    "Q7" = 0.1666666666666667
    "R7" = 0.05747126436781609
    "S7" = 0.396551724137931
    "B8" = 0.1222707423580786
    "D8" = 0.01528384279475982
    "F8" = 0.04803493449781659
    "J8" = 0.1135371179039301
    "O8" = 0.02838427947598253
    "Q8" = 0.1965065502183406
    "R8" = 0.07205240174672489
    "S8" = 0.4039301310043668
    "B9" = 0.1449275362318841
    "D9" = 0.01449275362318841
    "F9" = 0.06280193236714976
    "J9" = 0.1207729468599034
    "O9" = 0.02898550724637681
    "Q9" = 0.1932367149758454
    "R9" = 0.09178743961352658
    "S9" = 0.3429951690821256
    "B10" = 0.1239731142643764
    "D10" = 0.01941747572815534
    "F10" = 0.07617625093353249
    "J10" = 0.1097834204630321
    "O10" = 0.0186706497386109
    "Q10" = 0.2255414488424197
    "R10" = 0.06422703510082151
    "S10" = 0.3622106049290515
    "F11" = 0.004
    "G11" = 0.132
    "J11" = 0.096
    "K11" = 0.188
    "L11" = 0.5679999999999999
    "S11" = 0.012
    "G12" = 0.7876712328767124
    "J12" = 0.1712328767123288
    "K12" = 0.0136986301369863
    "L12" = 0.02054794520547945
    "S12" = 0.00684931506849315
    "G13" = 0.6956521739130435
    "J13" = 0.2826086956521739
    "S13" = 0.02173913043478261
    "F15" = 0.004081632653061225
    "H15" = 0.1224489795918367
    "I15" = 0.1061224489795918
    "J15" = 0.3836734693877551
    "K15" = 0.04081632653061224
    "M15" = 0.00816326530612245
    "O15" = 0.06122448979591837
    "S15" = 0.273469387755102
    "F16" = 0.02380952380952381
    "H16" = 0.2
    "I16" = 0.05714285714285714
    "J16" = 0.4619047619047619
    "K16" = 0.08571428571428572
    "M16" = 0.02857142857142857
    "O16" = 0.04761904761904762
    "S16" = 0.09523809523809523
    "F17" = 0.01181102362204724
    "H17" = 0.1830708661417323
    "I17" = 0.09251968503937008
    "J17" = 0.4330708661417323
    "K17" = 0.07874015748031496
    "M17" = 0.007874015748031496
    "O17" = 0.05511811023622047
    "S17" = 0.1377952755905512
    "F18" = 0.01257861635220126
    "H18" = 0.1320754716981132
    "I18" = 0.07547169811320754
    "J18" = 0.4339622641509434
    "K18" = 0.09433962264150944
    "O18" = 0.1069182389937107
    "S18" = 0.1446540880503145
    "F19" = 0.0076103500761035
    "H19" = 0.2100456621004566
    "I19" = 0.08447488584474885
    "J19" = 0.380517503805175
    "K19" = 0.0867579908675799
    "M19" = 0.02663622526636225
    "N19" = 0.0015220700152207
    "O19" = 0.06773211567732115
    "S19" = 0.134703196347032
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
